# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (col E) / "Valor Mora" (col F) table occupying rows
# 16-24 gets re-sorted from descending period order (2103..1901) to
# ascending period order (1901..2103). Columns B, C, D, G, H, I, J are
# identical across these rows, so the net effect is simply that the E/F
# pairs for the 9 rows are reversed top-to-bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 16
$lastRow  = 24

# Snapshot the current Periodo Mora (E) and Valor Mora (F) values before
# writing anything back, so the row-by-row reversal doesn't clobber
# values we still need to read.
$periodo = @{}
$valor   = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $periodo[$r] = $ws.Cells.Item($r, 5).Value2
    $valor[$r]   = $ws.Cells.Item($r, 6).Value2
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $firstRow + $lastRow - $r
    $ws.Cells.Item($r, 5).Value2 = $periodo[$srcRow]
    $ws.Cells.Item($r, 6).Value2 = $valor[$srcRow]
}
